$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos data (prices, volume %, and the Stellar/Monero row swap).
# Column D values that look like plain numbers must be forced to Text format
# first, otherwise Excel auto-converts them to numeric (losing trailing zeros /
# the original "number-as-text" representation used by the source data).

$ws.Range("D2").Value = "61.628.32"
$ws.Range("E2").Value = "  -3.40%  "

$ws.Range("D3").Value = "3.007.68"

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.36"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.40"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.001.89"
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("E10").Value = "  -5.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  -4.66%  "

$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.96"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "3.492.90"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").Value = "61.742.67"
$ws.Range("E16").Value = "  -3.31%  "

$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").Value = "3.008.87"
$ws.Range("E18").Value = "  -2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.24"
$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -3.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.01"
$ws.Range("E23").Value = "  -1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.45"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.04"
$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("E28").Value = "  -4.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.62"
$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.43"
$ws.Range("E34").Value = "  -2.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "458.23"
$ws.Range("E37").Value = "  -8.61%  "

$ws.Range("D38").Value = "3.199.98"
$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0793"
$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0382"
$ws.Range("E40").Value = "  -2.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.11"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -8.82%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.50"
$ws.Range("E45").Value = "  +4.63%  "

$ws.Range("E46").Value = "  -4.53%  "

$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.11"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "0.0₃0489"
$ws.Range("E50").Value = "  -7.66%  "

$ws.Range("E51").Value = "  +6.26%  "
